# Apply "Add data for 2022-07-19" update:
# - Rename the sheet and update the "through July NN" label from 10 to 11
# - Update the running-total counts for July 2022 (column B) and a handful
#   of scattered cells across several neighborhood rows to reflect the
#   newly-added day of carjacking data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet title / label for the partial "July 2022" month column.
$ws.Name = "Through 2022-07-11"
$ws.Range("B1").Value = "July 2022 (through July 11)"

# Row 2 - Austin
$ws.Range("B2").Value = 5
$ws.Range("AR2").Value = 3
$ws.Range("AY2").Value = 2

# Row 3 - Englewood
$ws.Range("B3").Value = 5
$ws.Range("AR3").Value = 3

# Row 4 - Auburn Gresham
$ws.Range("B4").Value = 3

# Row 6 - Grand Crossing
$ws.Range("B6").Value = 5
$ws.Range("W6").Value = 1

# Row 8 - North Lawndale
$ws.Range("B8").Value = 3
$ws.Range("P8").Value = 7
$ws.Range("AR8").Value = 2

# Row 11 - Loop
$ws.Range("P11").Value = 2
$ws.Range("W11").Value = 1

# Row 19 - South Shore
$ws.Range("I19").Value = 1
$ws.Range("AY19").Value = 1

# Row 23 - South Chicago
$ws.Range("P23").Value = 3

# Row 26 - Little Village
$ws.Range("B26").Value = 3

# Row 29 - Humboldt Park
$ws.Range("P29").Value = 1

# Row 31 - Fuller Park
$ws.Range("AK31").Value = 1

# Row 41 - Logan Square
$ws.Range("B41").Value = 1

# Row 44 - New City
$ws.Range("AR44").Value = 1

# Row 47 - Little Italy, UIC
$ws.Range("AD47").Value = 1

# Row 52 - Chatham
$ws.Range("B52").Value = 2

# Row 54 - Bucktown
$ws.Range("B54").Value = 1

# Row 71 - Galewood
$ws.Range("AK71").Value = 1

# Row 78 - Lake View
$ws.Range("B78").Value = 1

# Row 81 - Montclare
$ws.Range("AR81").Value = 1
